# Update gh-pages to output generated at 456a3b4
# Updates "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 420
$ws1.Range("F3").Value = 2655
$ws1.Range("F4").Value = 120

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 420
$ws4.Range("F7").Value = 2655
$ws4.Range("F8").Value = 120
